# Template cleanup: rename sheet, drop the leftover (empty, hyperlink-styled)
# rows below the header, resize the header columns, drop the unused 5th
# column width override, and move the saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "P" -> "Sheet1"
$ws.Name = "Sheet1"

# Rows 2:4 only ever held orphaned empty/hyperlink-styled cells (no data) -
# drop them so the sheet's used range shrinks back to the header row.
$ws.Rows("2:4").Delete()

# Re-size the header columns (A-D) and let the now-unused 5th column
# collapse back to the sheet's standard width.
$ws.Columns("A").ColumnWidth = 12.67
$ws.Columns("B").ColumnWidth = 15.67
$ws.Columns("C").ColumnWidth = 15.67
$ws.Columns("D").ColumnWidth = 10.67
$ws.Columns("E").ColumnWidth = $ws.StandardWidth

# Move the stored selection.
$null = $ws.Range("G7").Select()

# The deleted rows were the only cells using the built-in "Hyperlink" cell
# style (leftover formatting, never an actual hyperlink) - drop the now
# unused named style from the workbook.
$null = $wb.Styles.Item("Hyperlink").Delete()
